# Weekly update: insert a new price record for Berenjena (Terminal La Palmera
# de La Serena) as row 294, pushing the existing rows 294-316 down to 295-317.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 294, shifting rows 294:316 down to 295:317.
$ws.Rows.Item(294).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A294").Value = 8
$ws.Range("B294").Value = "Terminal La Palmera de La Serena"
$ws.Range("C294").Value = "Coquimbo"
$ws.Range("D294").Value = 45223
$ws.Range("E294").Value = 4
$ws.Range("F294").Value = 100112001
$ws.Range("G294").Value = "Berenjena"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 460
$ws.Range("K294").Value = 11000
$ws.Range("L294").Value = 12000
$ws.Range("M294").Value = 11500
$ws.Range("N294").Value = "$/caja 50 unidades"
$ws.Range("O294").Value = "Región de Arica y Parinacota"
$ws.Range("P294").Value = 230
$ws.Range("Q294").Value = 50
$ws.Range("R294").Value = "Hortaliza"

# Match the date-number-format style used by the rest of column D.
$ws.Range("D294").NumberFormat = $ws.Range("D295").NumberFormat
